{"js": "const replacements = [[\"2025-12-15 Monday\", \"2025-12-16 Tuesday\"], [\"398\u00d78=3184\", \"748\u00d76=4488\"], [\"332\u00d73=996\", \"302\u00d74=1208\"], [\"927\u00d77=6489\", \"660\u00d79=5940\"], [\"473\u00d73=1419\", \"714\u00d74=2856\"], [\"174\u00d79=1566\", \"127\u00d75=635\"], [\"533\u00d73=1599\", \"988\u00d74=3952\"], [\"564\u00d77=3948\", \"530\u00d78=4240\"], [\"684\u00d77=4788\", \"155\u00d79=1395\"], [\"288\u00d78=2304\", \"643\u00d78=5144\"], [\"168\u00d72=336\", \"796\u00d73=2388\"], [\"114\u00d74=456\", \"252\u00d76=1512\"], [\"665\u00d79=5985\", \"750\u00d78=6000\"], [\"709\u00d76=4254\", \"629\u00d75=3145\"], [\"261\u00d77=1827\", \"758\u00d72=1516\"], [\"211\u00d73=633\", \"640\u00d76=3840\"], [\"703\u00d78=5624\", \"674\u00d72=1348\"], [\"170\u00d76=1020\", \"209\u00d78=1672\"], [\"739\u00d75=3695\", \"518\u00d76=3108\"], [\"229\u00d74=916\", \"239\u00d79=2151\"], [\"272\u00d79=2448\", \"886\u00d73=2658\"], [\"879\u00d72=1758\", \"466\u00d72=932\"], [\"310\u00d73=930\", \"117\u00d72=234\"], [\"133\u00d73=399\", \"246\u00d74=984\"], [\"554\u00d77=3878\", \"321\u00d72=642\"], [\"968\u00d76=5808\", \"986\u00d78=7888\"]];\n\nconst body = context.document.body;\nlet totalReplaced = 0;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n    totalReplaced++;\n  }\n  await context.sync();\n}\n\nreturn \"replaced:\" + totalReplaced;\n", "ps1": "$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"2025-12-15 Monday\", \"2025-12-16 Tuesday\"),\n    @(\"398\u00d78=3184\", \"748\u00d76=4488\"),\n    @(\"332\u00d73=996\", \"302\u00d74=1208\"),\n    @(\"927\u00d77=6489\", \"660\u00d79=5940\"),\n    @(\"473\u00d73=1419\", \"714\u00d74=2856\"),\n    @(\"174\u00d79=1566\", \"127\u00d75=635\"),\n    @(\"533\u00d73=1599\", \"988\u00d74=3952\"),\n    @(\"564\u00d77=3948\", \"530\u00d78=4240\"),\n    @(\"684\u00d77=4788\", \"155\u00d79=1395\"),\n    @(\"288\u00d78=2304\", \"643\u00d78=5144\"),\n    @(\"168\u00d72=336\", \"796\u00d73=2388\"),\n    @(\"114\u00d74=456\", \"252\u00d76=1512\"),\n    @(\"665\u00d79=5985\", \"750\u00d78=6000\"),\n    @(\"709\u00d76=4254\", \"629\u00d75=3145\"),\n    @(\"261\u00d77=1827\", \"758\u00d72=1516\"),\n    @(\"211\u00d73=633\", \"640\u00d76=3840\"),\n    @(\"703\u00d78=5624\", \"674\u00d72=1348\"),\n    @(\"170\u00d76=1020\", \"209\u00d78=1672\"),\n    @(\"739\u00d75=3695\", \"518\u00d76=3108\"),\n    @(\"229\u00d74=916\", \"239\u00d79=2151\"),\n    @(\"272\u00d79=2448\", \"886\u00d73=2658\"),\n    @(\"879\u00d72=1758\", \"466\u00d72=932\"),\n    @(\"310\u00d73=930\", \"117\u00d72=234\"),\n    @(\"133\u00d73=399\", \"246\u00d74=984\"),\n    @(\"554\u00d77=3878\", \"321\u00d72=642\"),\n    @(\"968\u00d76=5808\", \"986\u00d78=7888\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Format = $false\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute($oldText, $false, $true, $false, $false, $false, $true, 1, $false, $newText, 2) | Out-Null\n}\n\nWrite-Output \"done\"\n"}
